$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.608.47'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '2.195.61'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.99'
$ws.Range("E5").Value = '  +5.57%  '

$ws.Range("E6").Value = '  +1.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.34'
$ws.Range("E7").Value = '  +4.60%  '

$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("E9").Value = '  +2.45%  '

$ws.Range("E10").Value = '  +1.47%  '

$ws.Range("E11").Value = '  +1.92%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.101'
$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.80'
$ws.Range("E13").Value = '  +1.84%  '

$ws.Range("D14").Value = '2.525.24'
$ws.Range("E14").Value = '  +1.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.30'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").Value = '2.193.26'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.780'
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("D18").Value = '42.481.62'
$ws.Range("E18").Value = '  +2.43%  '

$ws.Range("E19").Value = '  +1.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.15'
$ws.Range("E20").Value = '  +2.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.90'
$ws.Range("E21").Value = '  +2.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.74'
$ws.Range("E22").Value = '  +0.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.16'
$ws.Range("E23").Value = '  +9.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.39'
$ws.Range("E24").Value = '  -6.41%  '

$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.68'
$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.39'
$ws.Range("E27").Value = '  +3.22%  '

$ws.Range("E28").Value = '  +1.37%  '

$ws.Range("E29").Value = '  +2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.40'
$ws.Range("E30").Value = '  +12.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.39'
$ws.Range("E31").Value = '  -1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.08'
$ws.Range("E32").Value = '  +1.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0800'
$ws.Range("E33").Value = '  +3.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.16'
$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  +2.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.34'
$ws.Range("E37").Value = '  +1.92%  '

$ws.Range("E38").Value = '  +8.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.21'
$ws.Range("E39").Value = '  +0.79%  '

$ws.Range("E40").Value = '  +0.72%  '

$ws.Range("E41").Value = '  +5.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").Value = '  -1.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '58.90'
$ws.Range("E43").Value = '  +0.54%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.98'
$ws.Range("E44").Value = '  +7.74%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.46'
$ws.Range("E45").Value = '  +1.27%  '

$ws.Range("B46").Value = 'WOONetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.475'
$ws.Range("E46").Value = '  +20.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0977'
$ws.Range("E47").Value = '  +1.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  +10.82%  '

$ws.Range("E49").Value = '  +2.15%  '

$ws.Range("E50").Value = '  +1.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.06'
$ws.Range("E51").Value = '  +20.22%  '
